$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before C (Integral), which shifts Integral -> D and Time -> E
$ws.Columns("C").Insert()

# New "Variance" header and B^2 formulas in column C
$ws.Range("C1").Value2 = "Variance"
$ws.Range("C2").Formula = "=B2^2"
$ws.Range("C3:C11").Formula = "=B3^2"

# Row 13 "Avg": add Variance average, apply scientific number format to the whole row
$ws.Range("C13").Formula = "=AVERAGE(C2:C11)"
$ws.Range("B13:E13").NumberFormat = "0.00000E+00"

# Row 14 "STD": drop the STD-of-(B column) and STD-of-(old-Integral, now Variance) columns;
# keep the STD of Integral (now column D) and Time (now column E)
$ws.Range("B14").ClearContents()

# Row 15 "RMS": sqrt of the average variance
$ws.Range("A15").Value2 = "RMS"
$ws.Range("B15").Formula = "=SQRT(C13)"

# Column widths to roughly match the final layout
$ws.Columns("B").ColumnWidth = 11.666666666666666
$ws.Columns("C").ColumnWidth = 10.330729166666666
$ws.Columns("D").ColumnWidth = 11.498697916666666
$ws.Columns("E").ColumnWidth = 11.830729166666666

# Selection matches the final saved state
$ws.Range("B13:E13").Select() | Out-Null
